# Apply the workbook/table naming changes described by the commit:
#   - rename worksheet "lw_summary" -> "link_summary_sheet"
#   - rename worksheet "lw_links"   -> "link_sheet"
#   - rename the "links_table" label (sheet2!A1) -> "link_table"
#   - add a new bold cell style and apply it to a new (empty) cell E28
#     on the links sheet

$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item(1)
$wsLinks   = $wb.Worksheets.Item(2)

# 1) Rename the sheets
$wsSummary.Name = "link_summary_sheet"
$wsLinks.Name   = "link_sheet"

# 2) Rename the "links_table" label to "link_table" (cell A1 on the links sheet)
$wsLinks.Range("A1").Value = "link_table"

# 3) Add a new bold style and apply it to a new cell (E28) on the links sheet
$newCell = $wsLinks.Range("E28")
$newCell.Font.Bold = $true
